# The "Row" column (A) and "Is_New" column (B) were swapped throughout the
# sheet: what used to live in column A (the numeric row index) now lives in
# column B, and column A is left blank. The header labels on row 1 swap too.
#
# Additionally, five rows describing the (withdrawn) "US Core PMO
# ServiceRequest Profile" rows (old rows 95-99) are removed; the row that
# used to be row 100 (US Core Specimen Profile / Specimen.collection.collector)
# shifts up to become the new row 95 -- it already carried the post-swap
# A/B layout, so no further swap is needed for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 94   # rows 2..94 get a straightforward A/B swap
$lastCol = 44   # column AR

# --- Header row: swap the "Row" / "Is_New" labels ---
# A1 carries a leading BOM artifact (left over from the CSV import) that
# stays pinned to the A1 cell itself rather than travelling with either
# label's text, so the new labels are written out explicitly rather than
# via a naive value swap.
$bom = [char]0xFEFF
$ws.Cells.Item(1, 1).Value = "$bom" + "Is_New"
$ws.Cells.Item(1, 2).Value = "Row"

# --- Data rows 2..94: move A's value into B, blank out A ---
for ($r = 2; $r -le $lastRow; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $aVal
    $ws.Cells.Item($r, 1).Value = $null
}

# --- Remove the five obsolete "PMO ServiceRequest" rows (95-99) ---
# This shifts the former row 100 (Specimen row, already in the new A/B
# layout) up into row 95, matching the target sheet exactly.
$ws.Range("A95:AR99").EntireRow.Delete()
